$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.2011505038214305, 1.831162888379733, 5.606418433192214, 2.367787666407656, 2.38270351621109, 51),
    @(-0.1088417725019823, 1.857027656740065, 5.611761966788237, 2.368915778745255, 2.390439169655512, 50),
    @(0.2647668497699802, 1.70436734373276, 5.093964917022694, 2.256981372768215, 2.264625190576031, 49),
    @(0.02918659729334675, 1.699897465629806, 4.92493289209014, 2.219218982455346, 2.242509468478655, 48),
    @(0.2598537466771387, 1.778557183867834, 5.424263086057175, 2.329004741527414, 2.339484978954009, 47),
    @(0.0866664478830151, 1.819987475692422, 5.629107254141546, 2.372573972322369, 2.397190134413542, 46),
    @(0.31532128808607, 1.709693465785397, 4.913699006851471, 2.216686492684852, 2.218938204032733, 45),
    @(0.04114813685194464, 1.807961608101199, 5.315196701504934, 2.305471036795937, 2.331753284961301, 44),
    @(0.4060028782411715, 1.842179574410668, 5.639651418059248, 2.374795026535816, 2.367523130472429, 43),
    @(0.1000955268392197, 1.911941050400914, 5.487513789220115, 2.342544298240722, 2.368774388965271, 42)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}

$wb.Save()
